# Auto-generated cell updates derived from the authoritative diff.
# For each affected row, H/I/J/K/L/M/N market-price columns are refreshed;
# a few rows also gain or lose an M/N cell entirely (ClearContents removes the
# node outright, matching the diff's cell deletions).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3358.353
$ws.Range("I40").Value = 3549.2
$ws.Range("J40").Value = 3085.7144
$ws.Range("K40").Value = 3549.2
$ws.Range("L40").Value = 3085.7144
$ws.Range("M40").Value = -3374.2
$ws.Range("N40").Value = -3435.7144
$ws.Range("H75").Value = 20000
$ws.Range("J75").Value = 20000
$ws.Range("L75").Value = 20000
$ws.Range("N75").Value = -21872
$ws.Range("H78").Value = 20000
$ws.Range("J78").Value = 20000
$ws.Range("L78").Value = 60000
$ws.Range("N78").Value = -69360
$ws.Range("H98").Value = 1455.5
$ws.Range("I98").Value = 1455.5
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1455.5
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 42.5
$ws.Range("N98").ClearContents()
$ws.Range("H112").Value = 989084.2
$ws.Range("J112").Value = 1137391.9
$ws.Range("L112").Value = 3412175.7
$ws.Range("N112").Value = -3414391.7
$ws.Range("H122").Value = 1455.5
$ws.Range("I122").Value = 1455.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4366.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1916.5
$ws.Range("N122").ClearContents()
$ws.Range("H131").Value = 1360.9032
$ws.Range("I131").Value = 476
$ws.Range("J131").Value = 2000
$ws.Range("K131").Value = 1428
$ws.Range("L131").Value = 6000
$ws.Range("M131").Value = 3612
$ws.Range("N131").Value = -16080

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1267.375
$ws.Range("I2").Value = 1185.1666
$ws.Range("J2").Value = 1514
$ws.Range("K2").Value = 1185.1666
$ws.Range("L2").Value = 1514
$ws.Range("M2").Value = -1072.1666
$ws.Range("N2").Value = -1740
$ws.Range("H45").Value = 2306.0386
$ws.Range("I45").Value = 2315.5908
$ws.Range("J45").Value = 2253.5
$ws.Range("K45").Value = 2315.5908
$ws.Range("L45").Value = 2253.5
$ws.Range("M45").Value = -1938.5908
$ws.Range("N45").Value = -3007.5
$ws.Range("H61").Value = 1903.5172
$ws.Range("I61").Value = 1617.7646
$ws.Range("J61").Value = 2308.3333
$ws.Range("K61").Value = 1617.7646
$ws.Range("L61").Value = 2308.3333
$ws.Range("M61").Value = -1405.7646
$ws.Range("N61").Value = -2732.3333
$ws.Range("H63").Value = 3057.8572
$ws.Range("I63").Value = 2281
$ws.Range("K63").Value = 2281
$ws.Range("M63").Value = -1595
$ws.Range("H66").Value = 3057.8572
$ws.Range("I66").Value = 2281
$ws.Range("K66").Value = 11405
$ws.Range("M66").Value = -7973
$ws.Range("H97").Value = 1246.409
$ws.Range("I97").Value = 1462.6061
$ws.Range("J97").Value = 597.8182
$ws.Range("K97").Value = 1462.6061
$ws.Range("L97").Value = 597.8182
$ws.Range("M97").Value = -966.6061
$ws.Range("N97").Value = -1589.8182
$ws.Range("H98").Value = 17449.75
$ws.Range("J98").Value = 17449.75
$ws.Range("L98").Value = 17449.75
$ws.Range("N98").Value = -23439.75
$ws.Range("H101").Value = 19500
$ws.Range("J101").Value = 19500
$ws.Range("L101").Value = 19500
$ws.Range("N101").Value = -25990
$ws.Range("H110").Value = 1341.3334
$ws.Range("I110").Value = 1299.3334
$ws.Range("J110").Value = 1467.3334
$ws.Range("K110").Value = 1299.3334
$ws.Range("L110").Value = 1467.3334
$ws.Range("M110").Value = 745.6666
$ws.Range("N110").Value = -5557.3334
$ws.Range("H113").Value = 40149
$ws.Range("J113").Value = 40149
$ws.Range("L113").Value = 40149
$ws.Range("N113").Value = -48827
$ws.Range("H116").Value = 1267.375
$ws.Range("I116").Value = 1185.1666
$ws.Range("J116").Value = 1514
$ws.Range("K116").Value = 1185.1666
$ws.Range("L116").Value = 1514
$ws.Range("M116").Value = 1108.8334
$ws.Range("N116").Value = -6102
$ws.Range("H125").Value = 41723
$ws.Range("J125").Value = 41723
$ws.Range("L125").Value = 41723
$ws.Range("N125").Value = -51563
$ws.Range("H132").Value = 28569.281
$ws.Range("I132").Value = 31262.97
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 93788.91
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -91258.91
$ws.Range("N132").Value = -20058.5
$ws.Range("H136").Value = 1903.5172
$ws.Range("I136").Value = 1617.7646
$ws.Range("J136").Value = 2308.3333
$ws.Range("K136").Value = 4853.293799999999
$ws.Range("L136").Value = 6924.999899999999
$ws.Range("M136").Value = -2303.293799999999
$ws.Range("N136").Value = -12024.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1267.375
$ws.Range("I3").Value = 1185.1666
$ws.Range("J3").Value = 1514
$ws.Range("K3").Value = 1185.1666
$ws.Range("L3").Value = 1514
$ws.Range("M3").Value = -1071.1666
$ws.Range("N3").Value = -1742
$ws.Range("H107").Value = 1196.9269
$ws.Range("I107").Value = 861.3226
$ws.Range("J107").Value = 2237.3
$ws.Range("K107").Value = 861.3226
$ws.Range("L107").Value = 2237.3
$ws.Range("M107").Value = 1058.6774
$ws.Range("N107").Value = -6077.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 46260
$ws.Range("J52").Value = 46260
$ws.Range("L52").Value = 46260
$ws.Range("N52").Value = -46848
$ws.Range("H74").Value = 22778
$ws.Range("J74").Value = 22778
$ws.Range("L74").Value = 22778
$ws.Range("N74").Value = -24526
$ws.Range("H77").Value = 22778
$ws.Range("J77").Value = 22778
$ws.Range("L77").Value = 68334
$ws.Range("N77").Value = -77070
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H134").Value = 4881.579
$ws.Range("I134").Value = 5422.4
$ws.Range("J134").Value = 2853.5
$ws.Range("K134").Value = 16267.2
$ws.Range("L134").Value = 8560.5
$ws.Range("M134").Value = -13732.2
$ws.Range("N134").Value = -13630.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 212.5238
$ws.Range("I34").Value = 103.31579
$ws.Range("J34").Value = 1250
$ws.Range("K34").Value = 309.94737
$ws.Range("L34").Value = 3750
$ws.Range("M34").Value = -225.94737
$ws.Range("N34").Value = -3918
$ws.Range("H131").Value = 844.5909
$ws.Range("J131").Value = 1009.2069
$ws.Range("L131").Value = 3027.6207
$ws.Range("N131").Value = -13107.6207
$ws.Range("H132").Value = 1053.4
$ws.Range("I132").Value = 531.6667
$ws.Range("J132").Value = 1836
$ws.Range("K132").Value = 4785.0003
$ws.Range("L132").Value = 16524
$ws.Range("M132").Value = -2255.0003
$ws.Range("N132").Value = -21584
$ws.Range("H137").Value = 3046.125
$ws.Range("I137").Value = 1005.75
$ws.Range("J137").Value = 6446.75
$ws.Range("K137").Value = 3017.25
$ws.Range("L137").Value = 19340.25
$ws.Range("M137").Value = 2082.75
$ws.Range("N137").Value = -29540.25
$ws.Range("H140").Value = 1913.1818
$ws.Range("I140").Value = 1199.4445
$ws.Range("J140").Value = 5125
$ws.Range("K140").Value = 3598.3335
$ws.Range("L140").Value = 15375
$ws.Range("M140").Value = 1581.6665
$ws.Range("N140").Value = -25735
$ws.Range("H141").Value = 1714
$ws.Range("I141").Value = 642.5
$ws.Range("J141").Value = 6000
$ws.Range("K141").Value = 1927.5
$ws.Range("L141").Value = 18000
$ws.Range("M141").Value = 3252.5
$ws.Range("N141").Value = -28360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 187.63158
$ws.Range("I2").Value = 123
$ws.Range("J2").Value = 259.44446
$ws.Range("K2").Value = 123
$ws.Range("L2").Value = 259.44446
$ws.Range("M2").Value = -10
$ws.Range("N2").Value = -485.44446
$ws.Range("H93").Value = 7333.3335
$ws.Range("J93").Value = 7333.3335
$ws.Range("L93").Value = 7333.3335
$ws.Range("N93").Value = -11077.3335
$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -6340
$ws.Range("H122").Value = 1306.2142
$ws.Range("I122").Value = 1403.35
$ws.Range("J122").Value = 1063.375
$ws.Range("K122").Value = 4210.049999999999
$ws.Range("L122").Value = 3190.125
$ws.Range("M122").Value = -1760.049999999999
$ws.Range("N122").Value = -8090.125
$ws.Range("H132").Value = 3621.5386
$ws.Range("I132").Value = 2645.75
$ws.Range("J132").Value = 4055.2222
$ws.Range("K132").Value = 7937.25
$ws.Range("L132").Value = 12165.6666
$ws.Range("M132").Value = -5407.25
$ws.Range("N132").Value = -17225.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H68").Value = 16918996
$ws.Range("I68").Value = 37593490
$ws.Range("J68").Value = 3502.9092
$ws.Range("K68").Value = 37593490
$ws.Range("L68").Value = 3502.9092
$ws.Range("M68").Value = -37592741
$ws.Range("N68").Value = -5000.9092
$ws.Range("H71").Value = 16918996
$ws.Range("I71").Value = 37593490
$ws.Range("J71").Value = 3502.9092
$ws.Range("K71").Value = 187967450
$ws.Range("L71").Value = 17514.546
$ws.Range("M71").Value = -187963706
$ws.Range("N71").Value = -25002.546
$ws.Range("H76").Value = 3000
$ws.Range("I76").Value = 3000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -2662
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 3000
$ws.Range("I79").Value = 3000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -1830
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 30183.334
$ws.Range("J80").Value = 30183.334
$ws.Range("L80").Value = 30183.334
$ws.Range("N80").Value = -32179.334
$ws.Range("H83").Value = 30183.334
$ws.Range("J83").Value = 30183.334
$ws.Range("L83").Value = 90550.00199999999
$ws.Range("N83").Value = -100534.002

